# Update country statistics (new data pull + re-sort) and refresh timestamp
# per commit "Update countries & provincias Spain".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refresh "last updated" timestamp string
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Agosto de 2020 a las 23:37"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 4809490
$ws.Cells.Item(4, 3).Value = 45172
$ws.Cells.Item(4, 4).Value = 2375219
$ws.Cells.Item(4, 5).Value = 2275964
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 409
$ws.Cells.Item(4, 8).Value = 158307

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 2733677
$ws.Cells.Item(5, 3).Value = 24801
$ws.Cells.Item(5, 4).Value = 1884051
$ws.Cells.Item(5, 5).Value = 755522
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 488
$ws.Cells.Item(5, 8).Value = 94104

# Row 21: Alemania
$ws.Cells.Item(21, 1).Value = "Alemania"
$ws.Cells.Item(21, 2).Value = 211462
$ws.Cells.Item(21, 3).Value = 385
$ws.Cells.Item(21, 4).Value = 193600
$ws.Cells.Item(21, 5).Value = 8636
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 9226

# Row 29: Egipto
$ws.Cells.Item(29, 1).Value = "Egipto"
$ws.Cells.Item(29, 2).Value = 94483
$ws.Cells.Item(29, 3).Value = 167
$ws.Cells.Item(29, 4).Value = 42455
$ws.Cells.Item(29, 5).Value = 47163
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 31
$ws.Cells.Item(29, 8).Value = 4865

# Row 36: Israel
$ws.Cells.Item(36, 1).Value = "Israel"
$ws.Cells.Item(36, 2).Value = 72815
$ws.Cells.Item(36, 3).Value = 597
$ws.Cells.Item(36, 4).Value = 45677
$ws.Cells.Item(36, 5).Value = 26602
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 10
$ws.Cells.Item(36, 8).Value = 536

# Row 52: Barein
$ws.Cells.Item(52, 1).Value = "Barein"
$ws.Cells.Item(52, 2).Value = 41536
$ws.Cells.Item(52, 3).Value = 346
$ws.Cells.Item(52, 4).Value = 38666
$ws.Cells.Item(52, 5).Value = 2723
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 147

# Row 109: Zimbabue
$ws.Cells.Item(109, 1).Value = "Zimbabue"
$ws.Cells.Item(109, 2).Value = 3921
$ws.Cells.Item(109, 3).Value = 262
$ws.Cells.Item(109, 4).Value = 1016
$ws.Cells.Item(109, 5).Value = 2835
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 70

# Row 110: Libia
$ws.Cells.Item(110, 1).Value = "Libia"
$ws.Cells.Item(110, 2).Value = 3691
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 619
$ws.Cells.Item(110, 5).Value = 2992
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 80

# Row 111: Nicaragua
$ws.Cells.Item(111, 1).Value = "Nicaragua"
$ws.Cells.Item(111, 2).Value = 3672
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 2492
$ws.Cells.Item(111, 5).Value = 1064
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 116

# Row 121: Cabo Verde
$ws.Cells.Item(121, 1).Value = "Cabo Verde"
$ws.Cells.Item(121, 2).Value = 2547
$ws.Cells.Item(121, 3).Value = 67
$ws.Cells.Item(121, 4).Value = 1860
$ws.Cells.Item(121, 5).Value = 663
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 24

# Row 122: Mali
$ws.Cells.Item(122, 1).Value = "Mali"
$ws.Cells.Item(122, 2).Value = 2535
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 1941
$ws.Cells.Item(122, 5).Value = 470
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 124

# Row 129: Ruanda
$ws.Cells.Item(129, 1).Value = "Ruanda"
$ws.Cells.Item(129, 2).Value = 2062
$ws.Cells.Item(129, 3).Value = 20
$ws.Cells.Item(129, 4).Value = 1144
$ws.Cells.Item(129, 5).Value = 913
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 5

# Row 143: Angola
$ws.Cells.Item(143, 1).Value = "Angola"
$ws.Cells.Item(143, 2).Value = 1199
$ws.Cells.Item(143, 3).Value = 35
$ws.Cells.Item(143, 4).Value = 461
$ws.Cells.Item(143, 5).Value = 683
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 55

# Row 144: Uganda
$ws.Cells.Item(144, 1).Value = "Uganda"
$ws.Cells.Item(144, 2).Value = 1182
$ws.Cells.Item(144, 3).Value = 6
$ws.Cells.Item(144, 4).Value = 1045
$ws.Cells.Item(144, 5).Value = 133
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 4

# Row 145: Georgia
$ws.Cells.Item(145, 1).Value = "Georgia"
$ws.Cells.Item(145, 2).Value = 1177
$ws.Cells.Item(145, 3).Value = 6
$ws.Cells.Item(145, 4).Value = 955
$ws.Cells.Item(145, 5).Value = 205
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 17

# Row 149: Togo
$ws.Cells.Item(149, 1).Value = "Togo"
$ws.Cells.Item(149, 2).Value = 961
$ws.Cells.Item(149, 3).Value = 3
$ws.Cells.Item(149, 4).Value = 660
$ws.Cells.Item(149, 5).Value = 282
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 19
